$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    # Assign through a literal-text formula, then convert the formula
    # result back to a plain value via copy / paste-special. This avoids
    # Excel's automatic "looks like a number/currency" type coercion
    # (which would store the cell as a numeric value instead of text)
    # while also avoiding any residual NumberFormat/style changes that a
    # forced-text ( leading apostrophe ) entry would otherwise leave behind.
    $escaped = $text -replace '"', '""'
    $range.Formula = '="' + $escaped + '"'
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null
}

# Row 2: tracking number / expected rate / result updated for new Pre-Prod URL
Set-TextValue $ws.Range("P2") "320018569086"
Set-TextValue $ws.Range("Q2") "`$19.04"
Set-TextValue $ws.Range("R2") "PASS"

# Row 3: tracking number updated
Set-TextValue $ws.Range("P3") "320018590118"

# Row 5: tracking number / expected rate updated
Set-TextValue $ws.Range("P5") "320018567576"
Set-TextValue $ws.Range("Q5") "`$43.07"

$excel.CutCopyMode = 0
